# Update the "Förändrad" (Changed) date column (C) from 45175 (2023-09-06)
# to 45177 (2023-09-08) for every data row (rows 2 through 132).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 132) { $lastRow = 132 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45177
    }
}
